# Applies the cryptos-list refresh described in the commit:
#  - updates Price (D) / Volume 1h (E) figures for most rows
#  - swaps the NEARProtocol / ImmutableX rows (32 <-> 33)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new text would otherwise be auto-parsed as a Number by
# Excel (single "." decimal, no thousands separators). A leading apostrophe
# (doubled to '' inside a PowerShell single-quoted string) forces Excel to
# keep them as text, matching the original inlineStr cell type.

$ws.Range('D2').Value = '68.856.65'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '3.921.52'
$ws.Range('E3').Value = '  +4.02%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''605.33'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = '''165.39'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('D7').Value = '3.917.62'
$ws.Range('E7').Value = '  +3.98%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  -1.56%  '
$ws.Range('E10').Value = '  -4.36%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').Value = '''0.462'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').Value = '''37.23'
$ws.Range('E13').Value = '  -1.67%  '
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = '4.581.25'
$ws.Range('E15').Value = '  +4.25%  '
$ws.Range('D16').Value = '3.950.48'
$ws.Range('E16').Value = '  +5.02%  '
$ws.Range('D17').Value = '69.011.39'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').Value = '''17.14'
$ws.Range('E20').Value = '  -3.53%  '
$ws.Range('D21').Value = '''11.11'
$ws.Range('E21').Value = '  -2.29%  '
$ws.Range('D22').Value = '''488.33'
$ws.Range('E22').Value = '  -0.73%  '
$ws.Range('D23').Value = '''0.725'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').Value = '''0.0000165'
$ws.Range('E24').Value = '  +11.17%  '
$ws.Range('D25').Value = '''84.47'
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('D28').Value = '''10.14'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('D31').Value = '4.076.96'
$ws.Range('E31').Value = '  +4.36%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Value = '''7.89'
$ws.Range('E32').Value = '  -3.43%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '''2.39'
$ws.Range('E33').Value = '  -2.17%  '
$ws.Range('D34').Value = '''32.44'
$ws.Range('E34').Value = '  +1.65%  '
$ws.Range('D35').Value = '3.871.34'
$ws.Range('E35').Value = '  +4.60%  '
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('E37').Value = '  +2.44%  '
$ws.Range('E38').Value = '  +1.98%  '
$ws.Range('D39').Value = '''5.92'
$ws.Range('E39').Value = '  -0.78%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  -1.42%  '
$ws.Range('D42').Value = '''443.83'
$ws.Range('E42').Value = '  +3.81%  '
$ws.Range('E43').Value = '  -3.67%  '
$ws.Range('E44').Value = '  -0.64%  '
$ws.Range('D45').Value = '''48.46'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('D46').Value = '''8.47'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D48').Value = '''26.68'
$ws.Range('E48').Value = '  +13.45%  '
$ws.Range('D49').Value = '2.848.40'
$ws.Range('E49').Value = '  +1.76%  '
$ws.Range('D50').Value = '''141.93'
$ws.Range('E50').Value = '  +0.37%  '
$ws.Range('D51').Value = '''0.0358'
$ws.Range('E51').Value = '  +1.41%  '
